$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# --- Cell 1: "personne_age" -> "personne_naissance", split into
#     "p" / "ersonne_" / "naissance" runs (matching the target diff) ---
$cell1 = $tbl.Cell(7, 1)
$start1 = $cell1.Range.Start

# Replace "age" with "naissance"
$ageRange = $d.Range($start1 + 9, $start1 + 12)
$ageRange.Text = "naissance"

# Re-split "p" from "ersonne_naissance" (toggle formatting forces a run split
# without altering the visible formatting)
$restRange = $d.Range($start1 + 1, $start1 + 18)
$restRange.Font.Bold = 1
$restRange.Font.Bold = 0

# Split "naissance" from "ersonne_"
$naisRange = $d.Range($start1 + 9, $start1 + 18)
$naisRange.Font.Bold = 1
$naisRange.Font.Bold = 0

# --- Cell 2: "Age de la personne" -> "Date de naissance" + " de la personne" ---
$cell2 = $tbl.Cell(7, 2)
$start2 = $cell2.Range.Start

# Replace "Age" with "Date de naissance"
$headRange = $d.Range($start2, $start2 + 3)
$headRange.Text = "Date de naissance"

# Split off " de la personne" into its own run
$tailRange = $d.Range($start2 + 17, $start2 + 32)
$tailRange.Font.Bold = 1
$tailRange.Font.Bold = 0
